$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(800, 50, 0.5, 0.6, 4584, -100),
    @(800, 50, 0.5, 0.6, 2740, 0),
    @(800, 50, 0.5, 0.6, 2722, 0),
    @(800, 50, 0.5, 0.6, 4718, 0),
    @(800, 50, 0.5, 0.6, 4559, 0)
)

$startRow = 236
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
